$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1, Q1 with same formatting as the existing header cells ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows (2-25): swap I<->K and M<->O values, and fill new P, Q columns with 2 ---
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value()   # column I
    $kVal = $ws.Cells.Item($r, 11).Value()  # column K
    $ws.Cells.Item($r, 9).Value = $kVal
    $ws.Cells.Item($r, 11).Value = $iVal

    $mVal = $ws.Cells.Item($r, 13).Value()  # column M
    $oVal = $ws.Cells.Item($r, 15).Value()  # column O
    $ws.Cells.Item($r, 13).Value = $oVal
    $ws.Cells.Item($r, 15).Value = $mVal

    $ws.Cells.Item($r, 16).Value = 2      # column P
    $ws.Cells.Item($r, 17).Value = 2      # column Q
}
